$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = -5
$ws.Range("F7").Value = -6
$ws.Range("F8").Value = -3
$ws.Range("F11").Value = 5
$ws.Range("F16").Value = 6
$ws.Range("F18").Value = -5
$ws.Range("F28").Value = -1
$ws.Range("F30").Value = -1
$ws.Range("F33").Value = -1
